$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new changelog rows (10-14) with dates and notes.
# Use copy/paste-formats from an existing date cell (A6) so the date cells
# keep reusing the workbook's existing "short date" style instead of Excel
# minting a brand-new number format for the literal value assignment.
$ws.Range("A6").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A10").Value = 44788
$ws.Range("B10").Value = "Made a 3D-model for the LQFP-128 package (V9990)"

$ws.Range("A11").Value = 44789
$ws.Range("B11").Value = "Replaced some decoupling capacitors for bigger ones."

$ws.Range("B12").Value = "Changed 21MHz crystal to SMD."

$ws.Range("B13").Value = "Corrected and error with labeling two parts of a symbol different numbers."

$ws.Range("A14").Value = 44793
$ws.Range("B14").Value = "Started routing the PCB, adding cut-outs to separate the analog from digital signal-domains."

# Update the active selection to match the new editing position.
$ws.Range("B15").Select()
